# The "2024" tracker sheet keeps a rolling log of bank notification
# text/timestamp pairs per month. A new September entry was logged, which
# pushes the existing September_Details/September_Date (R:S) entries down
# by one row; the oldest September row in turn spills into what used to be
# the first August_Details/August_Date (P:Q) row, which itself cascades
# down through the August block, finally pushing the trailing "Broadband"
# label in column A down into a brand-new last row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Shift September_Details/September_Date (R35:S101) down by one row. Excel
# grows the used range (and the sheet's row count / dimension) to fit, so
# the R101:S101 entry spills into row 102, the old P102:Q102 ("hdfc" /
# "2024-08-30 12:15:48") cascades down the August column into row 103, and
# so on down to row 106, finally bumping the "Broadband" label from A106
# into the newly created A107.
$ws.Range("R35:S35").Insert()

# Record the brand-new September notification at the top of the list.
$ws.Range("R35").Value2 = "corporate internet share"
$ws.Range("S35").Value2 = "2024-09-09 11:32:23"
